$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename header cells: replace spaces with underscores
$ws.Range("A1").Value = "Player_Id"
$ws.Range("B1").Value = "Player_First_Name"
$ws.Range("C1").Value = "Player_Last_Name"

# Reset the view: scroll back to top-left and select C2
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C2").Select()
